$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 5998.75
$ws.Range("I43").Value = 5998.3335
$ws.Range("K43").Value = 5998.3335
$ws.Range("M43").Value = -5929.3335
$ws.Range("H112").Value = 2267.5
$ws.Range("I112").Value = 1990.5555
$ws.Range("J112").Value = 2766
$ws.Range("K112").Value = 5971.666499999999
$ws.Range("L112").Value = 8298
$ws.Range("M112").Value = -4863.666499999999
$ws.Range("N112").Value = -10514
$ws.Range("H132").Value = 10455.046
$ws.Range("I132").Value = 11364.3
$ws.Range("J132").Value = 1362.5
$ws.Range("K132").Value = 34092.89999999999
$ws.Range("L132").Value = 4087.5
$ws.Range("M132").Value = -31562.89999999999
$ws.Range("N132").Value = -9147.5
$ws.Range("H138").Value = 3851
$ws.Range("J138").Value = 3586.2334
$ws.Range("L138").Value = 10758.7002
$ws.Range("N138").Value = -21038.7002
$ws.Range("H141").Value = 6367.1
$ws.Range("I141").Value = 2667.2856
$ws.Range("K141").Value = 8001.8568
$ws.Range("M141").Value = -2821.8568

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 191
$ws.Range("I5").Value = 99
$ws.Range("J5").Value = 317.5
$ws.Range("K5").Value = 99
$ws.Range("L5").Value = 317.5
$ws.Range("M5").Value = 13
$ws.Range("N5").Value = -541.5
$ws.Range("H97").Value = 1129.4375
$ws.Range("I97").Value = 991.4
$ws.Range("K97").Value = 991.4
$ws.Range("M97").Value = -495.4
$ws.Range("H110").Value = 977.53845
$ws.Range("I110").Value = 977.53845
$ws.Range("K110").Value = 977.53845
$ws.Range("M110").Value = 1067.46155
$ws.Range("H122").Value = 2172.1875
$ws.Range("I122").Value = 1789.2693
$ws.Range("K122").Value = 5367.8079
$ws.Range("M122").Value = -2917.8079
$ws.Range("H132").Value = 5202.579
$ws.Range("I132").Value = 3517.423
$ws.Range("K132").Value = 10552.269
$ws.Range("M132").Value = -8022.269
$ws.Range("H139").Value = 68928.75
$ws.Range("J139").Value = 68928.75
$ws.Range("L139").Value = 68928.75
$ws.Range("N139").Value = -79208.75

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 191
$ws.Range("I4").Value = 99
$ws.Range("J4").Value = 317.5
$ws.Range("K4").Value = 99
$ws.Range("L4").Value = 317.5
$ws.Range("M4").Value = 16
$ws.Range("N4").Value = -547.5
$ws.Range("H22").Value = 400
$ws.Range("I22").Value = 400
$ws.Range("K22").Value = 400
$ws.Range("M22").Value = -227
$ws.Range("H96").Value = 39993.25
$ws.Range("I96").Value = 24986.5
$ws.Range("J96").Value = 55000
$ws.Range("K96").Value = 24986.5
$ws.Range("L96").Value = 55000
$ws.Range("M96").Value = -22240.5
$ws.Range("N96").Value = -60492
$ws.Range("H107").Value = 1777.75
$ws.Range("I107").Value = 1777.75
$ws.Range("K107").Value = 1777.75
$ws.Range("M107").Value = 142.25

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3841.5
$ws.Range("I31").Value = 3463.4546
$ws.Range("K31").Value = 3463.4546
$ws.Range("M31").Value = -3168.4546
$ws.Range("H34").Value = 3841.5
$ws.Range("I34").Value = 3463.4546
$ws.Range("K34").Value = 3463.4546
$ws.Range("M34").Value = -3261.4546
$ws.Range("H58").Value = 4988.457
$ws.Range("I58").Value = 4165.8335
$ws.Range("K58").Value = 4165.8335
$ws.Range("M58").Value = -3962.8335
$ws.Range("H86").Value = 118154.125
$ws.Range("J86").Value = 10648.333
$ws.Range("L86").Value = 10648.333
$ws.Range("N86").Value = -12894.333
$ws.Range("H89").Value = 118154.125
$ws.Range("J89").Value = 10648.333
$ws.Range("L89").Value = 53241.665
$ws.Range("N89").Value = -64473.665
$ws.Range("H94").Value = 10051.546
$ws.Range("I94").Value = 100000
$ws.Range("J94").Value = 1056.7
$ws.Range("K94").Value = 100000
$ws.Range("L94").Value = 1056.7
$ws.Range("M94").Value = -99549
$ws.Range("N94").Value = -1958.7
$ws.Range("H107").Value = 1624
$ws.Range("I107").Value = 1506.0834
$ws.Range("J107").Value = 1752.6364
$ws.Range("K107").Value = 1506.0834
$ws.Range("L107").Value = 1752.6364
$ws.Range("M107").Value = 413.9166
$ws.Range("N107").Value = -5592.6364
$ws.Range("H134").Value = 1848.7142
$ws.Range("I134").Value = 1342.8823
$ws.Range("K134").Value = 4028.6469
$ws.Range("M134").Value = -1493.6469
$ws.Range("H136").Value = 4988.457
$ws.Range("I136").Value = 4165.8335
$ws.Range("K136").Value = 12497.5005
$ws.Range("M136").Value = -9947.500499999998

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1160.4615
$ws.Range("J5").Value = 1917.5
$ws.Range("L5").Value = 5752.5
$ws.Range("N5").Value = -5976.5
$ws.Range("H37").Value = 41086.086
$ws.Range("J37").Value = 41086.086
$ws.Range("L37").Value = 123258.258
$ws.Range("N37").Value = -123482.258
$ws.Range("H135").Value = 1160.4615
$ws.Range("J135").Value = 1917.5
$ws.Range("L135").Value = 17257.5
$ws.Range("N135").Value = -22327.5
$ws.Range("H140").Value = 10000
$ws.Range("I140").Value = 10000
$ws.Range("K140").Value = 30000
$ws.Range("M140").Value = -24820

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H57").Value = 22124.5
$ws.Range("J57").Value = 26166
$ws.Range("L57").Value = 26166
$ws.Range("N57").Value = -27806
$ws.Range("H113").Value = 1891.7931
$ws.Range("I113").Value = 1892.25
$ws.Range("K113").Value = 1892.25
$ws.Range("M113").Value = 277.75
$ws.Range("H132").Value = 13034.5625
$ws.Range("I132").Value = 13132
$ws.Range("K132").Value = 39396
$ws.Range("M132").Value = -36866

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 1667816.9
$ws.Range("I19").Value = 1380.2
$ws.Range("K19").Value = 1380.2
$ws.Range("M19").Value = -1210.2
$ws.Range("H40").Value = 2500
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1864
$ws.Range("H46").Value = 3755.1333
$ws.Range("I46").Value = 1930.75
$ws.Range("J46").Value = 4971.3887
$ws.Range("K46").Value = 1930.75
$ws.Range("L46").Value = 4971.3887
$ws.Range("M46").Value = -1742.75
$ws.Range("N46").Value = -5347.3887
$ws.Range("H61").Value = 8924.4
$ws.Range("I61").Value = 9405.1875
$ws.Range("J61").Value = 7001.25
$ws.Range("K61").Value = 9405.1875
$ws.Range("L61").Value = 7001.25
$ws.Range("M61").Value = -9203.1875
$ws.Range("N61").Value = -7405.25
$ws.Range("H113").Value = 8924.4
$ws.Range("I113").Value = 9405.1875
$ws.Range("J113").Value = 7001.25
$ws.Range("K113").Value = 9405.1875
$ws.Range("L113").Value = 7001.25
$ws.Range("M113").Value = -7235.1875
$ws.Range("N113").Value = -11341.25
$ws.Range("H132").Value = 3468.5833
$ws.Range("I132").Value = 2451.1428
$ws.Range("J132").Value = 4893
$ws.Range("K132").Value = 7353.428400000001
$ws.Range("L132").Value = 14679
$ws.Range("M132").Value = -4823.428400000001
$ws.Range("N132").Value = -19739
$ws.Range("H136").Value = 9377.571
$ws.Range("I136").Value = 3664.4285
$ws.Range("J136").Value = 15090.714
$ws.Range("K136").Value = 10993.2855
$ws.Range("L136").Value = 45272.142
$ws.Range("M136").Value = -8443.2855
$ws.Range("N136").Value = -50372.142

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 55557.617
$ws.Range("I81").Value = 2867.8572
$ws.Range("K81").Value = 5735.7144
$ws.Range("M81").Value = -4674.7144
$ws.Range("H84").Value = 55557.617
$ws.Range("I84").Value = 2867.8572
$ws.Range("K84").Value = 28678.572
$ws.Range("M84").Value = -23374.572
$ws.Range("H107").Value = 1292.7906
$ws.Range("I107").Value = 994.48
$ws.Range("K107").Value = 2983.44
$ws.Range("M107").Value = -1063.44
$ws.Range("H132").Value = 4348.6665
$ws.Range("I132").Value = 3746
$ws.Range("J132").Value = 4650
$ws.Range("K132").Value = 11238
$ws.Range("L132").Value = 13950
$ws.Range("M132").Value = -8708
$ws.Range("N132").Value = -19010
